$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 205; this shifts the existing rows 205-220
# down to 206-221 and updates the sheet dimension automatically.
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new weekly record.
$ws.Range("A205").Value = 3
$ws.Range("B205").Value = 'Femacal de La Calera'
$ws.Range("C205").Value = 'Coquimbo'
$ws.Range("D205").Value = 44826
$ws.Range("E205").Value = 5
$ws.Range("F205").Value = 100112010
$ws.Range("G205").Value = 'Achicoria'
$ws.Range("H205").Value = 'Sin especificar'
$ws.Range("I205").Value = 'Primera'
$ws.Range("J205").Value = 115
$ws.Range("K205").Value = 6000
$ws.Range("L205").Value = 6500
$ws.Range("M205").Value = 6239
$ws.Range("N205").Value = '$/caja 16 unidades'
$ws.Range("O205").Value = 'Provincia de Quillota'
$ws.Range("P205").Value = 390
$ws.Range("Q205").Value = 16
$ws.Range("R205").Value = 'Hortaliza'
